$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 44379
$ws.Range("J2").Value = 120
$ws.Range("K2").Value = 12000
$ws.Range("L2").Value = 13000
$ws.Range("M2").Value = 12667
$ws.Range("P2").Value = 974

# Row 3
$ws.Range("D3").Value = 44320
$ws.Range("J3").Value = 160
$ws.Range("K3").Value = 19000
$ws.Range("L3").Value = 20000
$ws.Range("M3").Value = 19500
$ws.Range("P3").Value = 1500

# Row 4
$ws.Range("D4").Value = 44616
$ws.Range("J4").Value = 120
$ws.Range("K4").Value = 19000
$ws.Range("L4").Value = 20000
$ws.Range("M4").Value = 19500
$ws.Range("P4").Value = 1500

# Row 5
$ws.Range("D5").Value = 44389
$ws.Range("J5").Value = 120
$ws.Range("K5").Value = 12000
$ws.Range("L5").Value = 13000
$ws.Range("M5").Value = 12500
$ws.Range("P5").Value = 962

# Row 6
$ws.Range("D6").Value = 44764
$ws.Range("J6").Value = 200
$ws.Range("K6").Value = 12000
$ws.Range("L6").Value = 13000
$ws.Range("M6").Value = 12500
$ws.Range("P6").Value = 962

# Row 8
$ws.Range("D8").Value = 44159
$ws.Range("J8").Value = 100
$ws.Range("K8").Value = 23000
$ws.Range("L8").Value = 24000
$ws.Range("M8").Value = 23500
$ws.Range("P8").Value = 1808

# Row 9
$ws.Range("D9").Value = 44469
$ws.Range("J9").Value = 140
$ws.Range("K9").Value = 13000
$ws.Range("L9").Value = 14000
$ws.Range("M9").Value = 13500
$ws.Range("P9").Value = 1038

# Row 10
$ws.Range("D10").Value = 44893
$ws.Range("J10").Value = 900
$ws.Range("K10").Value = 13000
$ws.Range("L10").Value = 14000
$ws.Range("M10").Value = 13444
$ws.Range("P10").Value = 1034

# Row 11
$ws.Range("D11").Value = 44397
$ws.Range("J11").Value = 140
$ws.Range("K11").Value = 12500
$ws.Range("L11").Value = 13000
$ws.Range("M11").Value = 12750
$ws.Range("P11").Value = 981

# Row 12
$ws.Range("D12").Value = 44406
$ws.Range("J12").Value = 160
$ws.Range("K12").Value = 17000
$ws.Range("L12").Value = 18000
$ws.Range("M12").Value = 17500
$ws.Range("P12").Value = 1346

# Row 13
$ws.Range("D13").Value = 44592
$ws.Range("J13").Value = 120
$ws.Range("K13").Value = 12000
$ws.Range("L13").Value = 13000
$ws.Range("M13").Value = 12500
$ws.Range("P13").Value = 962

# Row 14
$ws.Range("D14").Value = 44832
$ws.Range("J14").Value = 100
$ws.Range("K14").Value = 13000
$ws.Range("L14").Value = 14000
$ws.Range("M14").Value = 13500
$ws.Range("P14").Value = 1038

# Row 15
$ws.Range("D15").Value = 44890
$ws.Range("J15").Value = 300
$ws.Range("K15").Value = 14000
$ws.Range("L15").Value = 15000
$ws.Range("M15").Value = 14500
$ws.Range("P15").Value = 1115

# Row 16
$ws.Range("D16").Value = 44580
$ws.Range("J16").Value = 160
$ws.Range("K16").Value = 11000
$ws.Range("L16").Value = 12000
$ws.Range("M16").Value = 11500
$ws.Range("P16").Value = 885
